$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.356.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.711.72"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.61"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5295"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06668"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.87"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07697"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.500"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.946.77"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.718.19"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5839"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8230"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.03"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.365.52"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "222.39"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.639"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.46"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.013"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.45"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.693"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.235"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.25"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05332"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.292"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.468"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.438"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.636"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.58%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9529"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.393"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5856"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.152.60"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +9.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01634"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.814"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8402"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.30"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.853.91"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.77"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4547"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.50%  "

# Rows 49 and 50 swap (EnergySwap <-> Frax) with updated values
$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.42%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.091"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.00%  "

# Row 51 updates
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05229"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.24%  "
